$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 501
$ws.Range("J12").Value = 501
$ws.Range("L12").Value = 501
$ws.Range("N12").Value = -841
$ws.Range("H82").Value = 1833.3334
$ws.Range("I82").Value = 1833.3334
$ws.Range("K82").Value = 5500.0002
$ws.Range("M82").Value = -5094.0002
$ws.Range("H85").Value = 1833.3334
$ws.Range("I85").Value = 1833.3334
$ws.Range("K85").Value = 5500.0002
$ws.Range("M85").Value = -4096.0002
$ws.Range("H116").Value = 4939.6
$ws.Range("I116").Value = 4900
$ws.Range("J116").Value = 4999
$ws.Range("K116").Value = 4900
$ws.Range("L116").Value = 4999
$ws.Range("M116").Value = -1458
$ws.Range("N116").Value = -11883
$ws.Range("H131").Value = 2270.3333
$ws.Range("I131").Value = 350
$ws.Range("K131").Value = 1050
$ws.Range("M131").Value = 3990
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 4852
$ws.Range("I137").Value = 4674.4
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 14023.2
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -11473.2
$ws.Range("N137").Value = -20100
$ws.Range("H138").Value = 2658.0833
$ws.Range("H141").Value = 2241.7778
$ws.Range("I141").Value = 2272.25
$ws.Range("K141").Value = 6816.75
$ws.Range("M141").Value = -1636.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1397.58
$ws.Range("I32").Value = 1039.1666
$ws.Range("K32").Value = 1039.1666
$ws.Range("M32").Value = -752.1666
$ws.Range("H61").Value = 3242.875
$ws.Range("I61").Value = 2989
$ws.Range("J61").Value = 3666
$ws.Range("K61").Value = 2989
$ws.Range("L61").Value = 3666
$ws.Range("M61").Value = -2777
$ws.Range("N61").Value = -4090
$ws.Range("H110").Value = 1643.4706
$ws.Range("I110").Value = 1495.6428
$ws.Range("J110").Value = 2333.3333
$ws.Range("K110").Value = 1495.6428
$ws.Range("L110").Value = 2333.3333
$ws.Range("M110").Value = 549.3571999999999
$ws.Range("N110").Value = -6423.3333
$ws.Range("H122").Value = 1137.3334
$ws.Range("I122").Value = 1137.3334
$ws.Range("K122").Value = 3412.0002
$ws.Range("M122").Value = -962.0001999999999
$ws.Range("H132").Value = 2037.2858
$ws.Range("I132").Value = 2002.2727
$ws.Range("K132").Value = 6006.8181
$ws.Range("M132").Value = -3476.8181
$ws.Range("H136").Value = 3242.875
$ws.Range("I136").Value = 2989
$ws.Range("J136").Value = 3666
$ws.Range("K136").Value = 8967
$ws.Range("L136").Value = 10998
$ws.Range("M136").Value = -6417
$ws.Range("N136").Value = -16098

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 183.28572
$ws.Range("I7").Value = 47.166668
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 47.166668
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 65.833332
$ws.Range("N7").Value = -1226
$ws.Range("H134").Value = 12753
$ws.Range("I134").Value = 12753
$ws.Range("K134").Value = 38259
$ws.Range("M134").Value = -35724

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 700
$ws.Range("I7").Value = 720.7692
$ws.Range("J7").Value = 565
$ws.Range("K7").Value = 720.7692
$ws.Range("L7").Value = 565
$ws.Range("M7").Value = -607.7692
$ws.Range("N7").Value = -791
$ws.Range("H22").Value = 725
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -1700
$ws.Range("H35").Value = 420000.66
$ws.Range("I35").Value = 558334.25
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 558334.25
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -558040.25
$ws.Range("N35").Value = -5588
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 537
$ws.Range("I105").Value = 471.8
$ws.Range("J105").Value = 700
$ws.Range("K105").Value = 471.8
$ws.Range("L105").Value = 700
$ws.Range("M105").Value = 1275.2
$ws.Range("N105").Value = -4194
$ws.Range("H122").Value = 1605.7142
$ws.Range("I122").Value = 1330.3334
$ws.Range("K122").Value = 3991.0002
$ws.Range("M122").Value = -1541.0002
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 8869.5
$ws.Range("I134").Value = 8869.5
$ws.Range("K134").Value = 26608.5
$ws.Range("M134").Value = -24073.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1038.8
$ws.Range("I50").Value = 231.33333
$ws.Range("J50").Value = 2250
$ws.Range("K50").Value = 693.99999
$ws.Range("L50").Value = 6750
$ws.Range("M50").Value = -212.99999
$ws.Range("N50").Value = -7712
$ws.Range("H53").Value = 1038.8
$ws.Range("I53").Value = 231.33333
$ws.Range("J53").Value = 2250
$ws.Range("K53").Value = 693.99999
$ws.Range("L53").Value = 6750
$ws.Range("M53").Value = -212.99999
$ws.Range("N53").Value = -7712
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H104").Value = 1400
$ws.Range("I104").Value = 1400
$ws.Range("K104").Value = 4200
$ws.Range("M104").Value = -1579
$ws.Range("H121").Value = 900
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 900
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2700
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5320
$ws.Range("H131").Value = 1413.8334
$ws.Range("J131").Value = 1995
$ws.Range("L131").Value = 5985
$ws.Range("N131").Value = -16065

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 340.4
$ws.Range("I2").Value = 457.35715
$ws.Range("K2").Value = 457.35715
$ws.Range("M2").Value = -344.35715
$ws.Range("H59").Value = 11000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 11000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 11000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -12166
$ws.Range("H97").Value = 6333
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 6333
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6333
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -7325
$ws.Range("H126").Value = 7841.25
$ws.Range("I126").Value = 4161.875
$ws.Range("K126").Value = 12485.625
$ws.Range("M126").Value = -10015.625
$ws.Range("H132").Value = 2081.9167
$ws.Range("I132").Value = 2100.4
$ws.Range("J132").Value = 1989.5
$ws.Range("K132").Value = 6301.200000000001
$ws.Range("L132").Value = 5968.5
$ws.Range("M132").Value = -3771.200000000001
$ws.Range("N132").Value = -11028.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3396.2856
$ws.Range("I16").Value = 3194.4
$ws.Range("J16").Value = 3901
$ws.Range("K16").Value = 3194.4
$ws.Range("L16").Value = 3901
$ws.Range("M16").Value = -3024.4
$ws.Range("N16").Value = -4241
$ws.Range("H21").Value = 41000
$ws.Range("J21").Value = 41000
$ws.Range("L21").Value = 41000
$ws.Range("N21").Value = -41348
$ws.Range("H30").Value = 741.25
$ws.Range("I30").Value = 388.33334
$ws.Range("J30").Value = 1800
$ws.Range("K30").Value = 388.33334
$ws.Range("L30").Value = 1800
$ws.Range("M30").Value = -280.33334
$ws.Range("N30").Value = -2016
$ws.Range("H35").Value = 201115.8
$ws.Range("I35").Value = 250519.75
$ws.Range("J35").Value = 3500
$ws.Range("K35").Value = 250519.75
$ws.Range("L35").Value = 3500
$ws.Range("M35").Value = -250183.75
$ws.Range("N35").Value = -4172
$ws.Range("H46").Value = 2932.35
$ws.Range("I46").Value = 2768.0908
$ws.Range("J46").Value = 3133.111
$ws.Range("K46").Value = 2768.0908
$ws.Range("L46").Value = 3133.111
$ws.Range("M46").Value = -2580.0908
$ws.Range("N46").Value = -3509.111
$ws.Range("H61").Value = 4278.6
$ws.Range("I61").Value = 3798
$ws.Range("K61").Value = 3798
$ws.Range("M61").Value = -3596
$ws.Range("H82").Value = 1010.6667
$ws.Range("I82").Value = 1010.6667
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1010.6667
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -649.6667
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1010.6667
$ws.Range("I85").Value = 1010.6667
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1010.6667
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 237.3333
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 4278.6
$ws.Range("I113").Value = 3798
$ws.Range("K113").Value = 3798
$ws.Range("M113").Value = -1628
$ws.Range("H136").Value = 2982.2
$ws.Range("I136").Value = 2982.2
$ws.Range("K136").Value = 8946.599999999999
$ws.Range("M136").Value = -6396.599999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 616.0833
$ws.Range("I100").Value = 414.14285
$ws.Range("J100").Value = 898.8
$ws.Range("K100").Value = 828.2857
$ws.Range("L100").Value = 1797.6
$ws.Range("M100").Value = -287.2857
$ws.Range("N100").Value = -2879.6
